$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (C) and P_Value (D) values per corrected Diebold-Mariano test
$ws.Range("C2").Value = -1.041823831080326
$ws.Range("D2").Value = 0.3088075895544602

$ws.Range("C3").Value = -0.3506101329108211
$ws.Range("D3").Value = 0.7292157686421357

$ws.Range("C4").Value = -0.2007990298098971
$ws.Range("D4").Value = 0.8427005179111091

$ws.Range("C5").Value = -0.1164142192603153
$ws.Range("D5").Value = 0.9083804309342023

$ws.Range("C6").Value = 0.8054786037009706
$ws.Range("D6").Value = 0.4291639593471186

$ws.Range("C7").Value = 1.090657346350184
$ws.Range("D7").Value = 0.2872255518190878

$ws.Range("C8").Value = 1.305853788508668
$ws.Range("D8").Value = 0.205092775619633

$ws.Range("C9").Value = 0.1369906012062921
$ws.Range("D9").Value = 0.8922841831313122

$ws.Range("C10").Value = 0.1909808744941197
$ws.Range("D10").Value = 0.8502920777729419

$ws.Range("C11").Value = 0.06135251729998123
$ws.Range("D11").Value = 0.9516324346901412
